$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 7999
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H135").Value = 13515869
$ws.Range("I135").Value = 18520878
$ws.Range("J135").Value = 2346.1
$ws.Range("K135").Value = 166687902
$ws.Range("L135").Value = 21114.9
$ws.Range("M135").Value = -166685367
$ws.Range("N135").Value = -26184.9
$ws.Range("H137").Value = 9101504
$ws.Range("I137").Value = 11776696
$ws.Range("J137").Value = 5849.6
$ws.Range("K137").Value = 35330088
$ws.Range("L137").Value = 17548.8
$ws.Range("M137").Value = -35327538
$ws.Range("N137").Value = -22648.8
$ws.Range("H138").Value = 8778.888999999999
$ws.Range("I138").Value = 6346.1665
$ws.Range("J138").Value = 9265.433999999999
$ws.Range("K138").Value = 19038.4995
$ws.Range("L138").Value = 27796.302
$ws.Range("M138").Value = -13898.4995
$ws.Range("N138").Value = -38076.302

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4637.1333
$ws.Range("I2").Value = 4606.7
$ws.Range("J2").Value = 4698
$ws.Range("K2").Value = 4606.7
$ws.Range("L2").Value = 4698
$ws.Range("M2").Value = -4493.7
$ws.Range("N2").Value = -4924
$ws.Range("H45").Value = 5323.0713
$ws.Range("I45").Value = 3503.8333
$ws.Range("J45").Value = 6687.5
$ws.Range("K45").Value = 3503.8333
$ws.Range("L45").Value = 6687.5
$ws.Range("M45").Value = -3126.8333
$ws.Range("N45").Value = -7441.5
$ws.Range("H61").Value = 2469.0908
$ws.Range("I61").Value = 2025.1724
$ws.Range("J61").Value = 5687.5
$ws.Range("K61").Value = 2025.1724
$ws.Range("L61").Value = 5687.5
$ws.Range("M61").Value = -1813.1724
$ws.Range("N61").Value = -6111.5
$ws.Range("H74").Value = 1561.9286
$ws.Range("I74").Value = 1013.4091
$ws.Range("J74").Value = 3573.1667
$ws.Range("K74").Value = 1013.4091
$ws.Range("L74").Value = 3573.1667
$ws.Range("M74").Value = -139.4091
$ws.Range("N74").Value = -5321.1667
$ws.Range("H77").Value = 1561.9286
$ws.Range("I77").Value = 1013.4091
$ws.Range("J77").Value = 3573.1667
$ws.Range("K77").Value = 5067.0455
$ws.Range("L77").Value = 17865.8335
$ws.Range("M77").Value = -699.0455000000002
$ws.Range("N77").Value = -26601.8335
$ws.Range("H116").Value = 4637.1333
$ws.Range("I116").Value = 4606.7
$ws.Range("J116").Value = 4698
$ws.Range("K116").Value = 4606.7
$ws.Range("L116").Value = 4698
$ws.Range("M116").Value = -2312.7
$ws.Range("N116").Value = -9286
$ws.Range("H136").Value = 2469.0908
$ws.Range("I136").Value = 2025.1724
$ws.Range("J136").Value = 5687.5
$ws.Range("K136").Value = 6075.5172
$ws.Range("L136").Value = 17062.5
$ws.Range("M136").Value = -3525.5172
$ws.Range("N136").Value = -22162.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4637.1333
$ws.Range("I3").Value = 4606.7
$ws.Range("J3").Value = 4698
$ws.Range("K3").Value = 4606.7
$ws.Range("L3").Value = 4698
$ws.Range("M3").Value = -4492.7
$ws.Range("N3").Value = -4926
$ws.Range("H86").Value = 69763.92999999999
$ws.Range("I86").Value = 85206.836
$ws.Range("J86").Value = 7992.3335
$ws.Range("K86").Value = 85206.836
$ws.Range("L86").Value = 7992.3335
$ws.Range("M86").Value = -84083.836
$ws.Range("N86").Value = -10238.3335
$ws.Range("H89").Value = 69763.92999999999
$ws.Range("I89").Value = 85206.836
$ws.Range("J89").Value = 7992.3335
$ws.Range("K89").Value = 426034.18
$ws.Range("L89").Value = 39961.6675
$ws.Range("M89").Value = -420418.18
$ws.Range("N89").Value = -51193.6675
$ws.Range("H134").Value = 2561
$ws.Range("I134").Value = 2339.5117
$ws.Range("J134").Value = 2957.8333
$ws.Range("K134").Value = 7018.5351
$ws.Range("L134").Value = 8873.499899999999
$ws.Range("M134").Value = -4483.5351
$ws.Range("N134").Value = -13943.4999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2759.5095
$ws.Range("I31").Value = 2539
$ws.Range("J31").Value = 2956.3928
$ws.Range("K31").Value = 2539
$ws.Range("L31").Value = 2956.3928
$ws.Range("M31").Value = -2244
$ws.Range("N31").Value = -3546.3928
$ws.Range("H34").Value = 2759.5095
$ws.Range("I34").Value = 2539
$ws.Range("J34").Value = 2956.3928
$ws.Range("K34").Value = 2539
$ws.Range("L34").Value = 2956.3928
$ws.Range("M34").Value = -2337
$ws.Range("N34").Value = -3360.3928
$ws.Range("H58").Value = 2102.7778
$ws.Range("I58").Value = 1999.7368
$ws.Range("J58").Value = 2347.5
$ws.Range("K58").Value = 1999.7368
$ws.Range("L58").Value = 2347.5
$ws.Range("M58").Value = -1796.7368
$ws.Range("N58").Value = -2753.5
$ws.Range("H134").Value = 2149.3691
$ws.Range("I134").Value = 2244.0176
$ws.Range("K134").Value = 6732.0528
$ws.Range("M134").Value = -4197.0528
$ws.Range("H136").Value = 2102.7778
$ws.Range("I136").Value = 1999.7368
$ws.Range("J136").Value = 2347.5
$ws.Range("K136").Value = 5999.2104
$ws.Range("L136").Value = 7042.5
$ws.Range("M136").Value = -3449.2104
$ws.Range("N136").Value = -12142.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1000.0625
$ws.Range("I5").Value = 808.8
$ws.Range("J5").Value = 1318.8334
$ws.Range("K5").Value = 2426.4
$ws.Range("L5").Value = 3956.5002
$ws.Range("M5").Value = -2314.4
$ws.Range("N5").Value = -4180.5002
$ws.Range("H31").Value = 94
$ws.Range("I31").Value = 94
$ws.Range("K31").Value = 282
$ws.Range("M31").Value = 6
$ws.Range("H68").Value = 1752.8889
$ws.Range("J68").Value = 2129.75
$ws.Range("L68").Value = 6389.25
$ws.Range("N68").Value = -8011.25
$ws.Range("H71").Value = 1752.8889
$ws.Range("J71").Value = 2129.75
$ws.Range("L71").Value = 19167.75
$ws.Range("N71").Value = -27279.75
$ws.Range("H135").Value = 1000.0625
$ws.Range("I135").Value = 808.8
$ws.Range("J135").Value = 1318.8334
$ws.Range("K135").Value = 7279.2
$ws.Range("L135").Value = 11869.5006
$ws.Range("M135").Value = -4744.2
$ws.Range("N135").Value = -16939.5006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9061.6
$ws.Range("I132").Value = 8407
$ws.Range("J132").Value = 12498.25
$ws.Range("K132").Value = 25221
$ws.Range("L132").Value = 37494.75
$ws.Range("M132").Value = -22691
$ws.Range("N132").Value = -42554.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4473.537
$ws.Range("I132").Value = 4385.5527
$ws.Range("J132").Value = 4682.5
$ws.Range("K132").Value = 13156.6581
$ws.Range("L132").Value = 14047.5
$ws.Range("M132").Value = -10626.6581
$ws.Range("N132").Value = -19107.5
$ws.Range("H136").Value = 5420.4614
$ws.Range("I136").Value = 4000
$ws.Range("J136").Value = 5678.727
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 17036.181
$ws.Range("M136").Value = -9450
$ws.Range("N136").Value = -22136.181

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3365.3667
$ws.Range("I132").Value = 3011.125
$ws.Range("J132").Value = 4782.3335
$ws.Range("K132").Value = 9033.375
$ws.Range("L132").Value = 14347.0005
$ws.Range("M132").Value = -6503.375
$ws.Range("N132").Value = -19407.0005
$ws.Range("H136").Value = 5416.8945
$ws.Range("I136").Value = 5181.3667
$ws.Range("J136").Value = 6300.125
$ws.Range("K136").Value = 15544.1001
$ws.Range("L136").Value = 18900.375
$ws.Range("M136").Value = -12994.1001
$ws.Range("N136").Value = -24000.375

Write-Output "Applied 195 cell updates across 8 sheets"